$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the AutoFilter from showing MYT1L rows to showing SHANK3B rows.
# Using the "values" filter mode (7 = xlFilterValues) reproduces the
# <filters><filter val="..."/></filters> shape Excel writes, and updates
# each row's Hidden state for us (rows 2-21 become hidden, rows 62-81
# become visible).
[void]$ws.Range("A1:L162").AutoFilter(2, @("SHANK3B"), 7)

# The previously-filtered rows (62-81, Gene = SHANK3B) had a one-off fill
# style applied to A63:A67; clear it back to the default "Normal" style
# now that the rows are visible again.
$ws.Range("A63:A67").Style = "Normal"

Write-Host "done"
